$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.6468346118927
$ws.Range("B1").Value = 2.909249305725098
$ws.Range("C1").Value = 4.695108413696289
$ws.Range("D1").Value = 1.401321053504944
$ws.Range("E1").Value = 0.8185690641403198
